# Commit: "added sample test to retrieve entire sheet"
#
# - Rename the existing sheet (Sheet1) to "TestData"
# - Fix the typo in cell A1 from "Testcases" to "TestCases"
# - Add two new, empty worksheets: "Sample" and "Demo" (after TestData)
# - Move the selection on TestData from B6 to B7
# - Keep TestData as the active/selected sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename Sheet1 -> TestData
$ws.Name = "TestData"

# Fix header typo: "Testcases" -> "TestCases"
$ws.Range("A1").Value = "TestCases"

# Add the two new sheets right after TestData, in order: Sample, then Demo
$sample = $wb.Worksheets.Add($null, $ws)
$sample.Name = "Sample"

$demo = $wb.Worksheets.Add($null, $sample)
$demo.Name = "Demo"

# Re-activate TestData and move its selection from B6 to B7
$ws.Activate()
$ws.Range("B7").Select()
